# edit.ps1
# Applies the "Add files via upload" header edit to the default header
# (header2.xml / Sections(1).Headers(wdHeaderFooterPrimary)):
#   1. Collapses the "Don't" / " Panic Room" / <w:br/> run trio (which were
#      wrapped in <w:proofErr> spell-check tags) into a single run reading
#      "Don't Panic Room" (no more manual line break, no more proofErr).
#   2. Replaces the stray combining-macron "char  KI" + " in der Kita" runs
#      (also wrapped in <w:proofErr> grammar tags) with a clean en dash
#      " - KI in der Kita" split across four runs (no more proofErr).
#   3. Tags every run (and each paragraph mark) in the header with
#      <w:lang w:val="en-GB"/>.
#   4. Appends a brand-new second paragraph to the header reading
#      "02Station Prompt the Prompt", split across two runs ("02S" /
#      "tation Prompt the Prompt"), matching the style of paragraph 1.
#
# The picture anchored in the header (the "KI-generierte Inhalte" image)
# is left completely untouched.
#
# Implementation note: this host's Range.LanguageID setter only stamps the
# <w:lang> element onto runs that already contain text -- it never reaches
# the paragraph-mark run properties stored in <w:pPr><w:rPr>. The only way
# found to reliably reproduce the paragraph-mark <w:lang> (and to drop the
# <w:proofErr/> bookmarks and <w:br/> cleanly) is to rebuild the header's
# two paragraphs verbatim via Range.InsertXML with a full WordProcessingML
# package fragment, reusing the existing picture markup unchanged.

$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1
$hdr = $d.Sections(1).Headers(1)
$targetRange = $hdr.Range

$newHeaderXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" mc:Ignorable="w14 w15 wp14"><w:body><w:p w14:paraId="1470AF61" w14:textId="2D9DFAC2" w:rsidR="0063194B" w:rsidRPr="00F017E9" w:rsidRDefault="00F017E9" w:rsidP="00F017E9"><w:pPr><w:pStyle w:val="Kopfzeile"/><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00F017E9"><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t>Don’t Panic Room</w:t></w:r><w:r w:rsidR="0063194B" w:rsidRPr="00F017E9"><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="331EC226" wp14:editId="7010C1B5"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>5655652</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>-168812</wp:posOffset></wp:positionV><wp:extent cx="611945" cy="611945"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="478329058" name="Grafik 478329058" descr="Ein Bild, das Astronomisches Objekt, Kugel, Planet, Astronomisches Ereignis enthält.&#xA;&#xA;KI-generierte Inhalte können fehlerhaft sein."/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="3" name="Grafik 3" descr="Ein Bild, das Astronomisches Objekt, Kugel, Planet, Astronomisches Ereignis enthält.&#xA;&#xA;KI-generierte Inhalte können fehlerhaft sein."/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="611945" cy="611945"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:anchor></w:drawing>
</w:r><w:r w:rsidRPr="00F017E9"><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r w:rsidRPr="00F017E9"><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t>KI</w:t></w:r><w:r w:rsidRPr="00F017E9"><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00F017E9"><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t>in der Kita</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kopfzeile"/><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t>02S</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-GB"/></w:rPr><w:t>tation Prompt the Prompt</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($newHeaderXml)

Write-Host "Header updated. Paragraph count:" $hdr.Range.Paragraphs.Count
Write-Host "Header text:" $hdr.Range.Text
